# BreakoutRoom2Instructions.pptx
# Slide 2, "Content Placeholder 2": the first bullet line currently reads
#   "You are facilitating a breakout room with one of the six sw quality dimensions"
# (split across 3 runs, with "sw" carrying a spell-check err="1" flag because it's
# an abbreviation). Replace it with the spelled-out
#   "You are facilitating a breakout room with one of the six Research Software (RS)  quality dimensions"
# as a single clean run (matching what PowerPoint produces when a user retypes the
# whole line), leaving every other paragraph/run in the text box untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$newFirstLine = "You are facilitating a breakout room with one of the six Research Software (RS)  quality dimensions"

# Replace paragraph 1 wholesale: delete it (runs + its trailing paragraph mark)
# and insert the replacement text + a paragraph mark before what is now the new
# paragraph 1 ("Create a google like the one on ..."). Doing it this way (rather
# than assigning .Text on the existing paragraph range) collapses the three
# original runs into a single new run, same as the authored edit.
$para1 = $tr.Paragraphs(1, 1)
[void]$para1.Delete()
[void]$tr.InsertBefore($newFirstLine + "`r")
